# Sprint3 Final Presentation was "shortened drastically": remove the
# duplicate "Project Goals" slide, the blank placeholder slide, and the
# four "Sprint 3 - Description of tasks" slides, leaving:
#   1 PokerGame
#   2 Problem Motivation
#   3 Project Goals
#   4 Client-Server Model      (was 5)
#   5 (full-slide picture)     (was 6)
#   6 DEMO                     (was 7)
#   7 Thank you!                (was 13)

$p = $ppt.ActivePresentation

# Delete from the highest index down so earlier indices stay valid.
$p.Slides.Item(12).Delete()  # Sprint 3 - Description of tasks (4/4)
$p.Slides.Item(11).Delete()  # Sprint 3 - Description of tasks (3/4)
$p.Slides.Item(10).Delete()  # Sprint 3 - Description of tasks (2/4)
$p.Slides.Item(9).Delete()   # Sprint 3 - Description of tasks (1/4)
$p.Slides.Item(8).Delete()   # blank placeholder slide
$p.Slides.Item(4).Delete()   # duplicate "Project Goals" slide
